$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1173
$ws.Range("I92").Value = 1173
$ws.Range("K92").Value = 1173
$ws.Range("M92").Value = 75

$ws.Range("H137").Value = 5040.3335
$ws.Range("I137").Value = 3119.4285
$ws.Range("K137").Value = 9358.2855
$ws.Range("M137").Value = -6808.2855

$ws.Range("H138").Value = 1668417.8
$ws.Range("I138").Value = 843.8333
$ws.Range("J138").Value = 4169778.5
$ws.Range("K138").Value = 2531.4999
$ws.Range("L138").Value = 12509335.5
$ws.Range("M138").Value = 2608.5001
$ws.Range("N138").Value = -12519615.5


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11115833
$ws.Range("I32").Value = 11115833
$ws.Range("K32").Value = 11115833
$ws.Range("M32").Value = -11115546

$ws.Range("H61").Value = 31317444
$ws.Range("I61").Value = 55557224
$ws.Range("K61").Value = 55557224
$ws.Range("M61").Value = -55557012

$ws.Range("H132").Value = 7257.1304
$ws.Range("I132").Value = 4124.9443
$ws.Range("K132").Value = 12374.8329
$ws.Range("M132").Value = -9844.832900000001

$ws.Range("H136").Value = 31317444
$ws.Range("I136").Value = 55557224
$ws.Range("K136").Value = 166671672
$ws.Range("M136").Value = -166669122


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 123000
$ws.Range("J62").Value = 126000
$ws.Range("L62").Value = 126000
$ws.Range("N62").Value = -127372

$ws.Range("H63").Value = 96000
$ws.Range("I63").Value = 70000
$ws.Range("J63").Value = 122000
$ws.Range("K63").Value = 70000
$ws.Range("L63").Value = 122000
$ws.Range("M63").Value = -69314
$ws.Range("N63").Value = -123372

$ws.Range("H65").Value = 123000
$ws.Range("J65").Value = 126000
$ws.Range("L65").Value = 378000
$ws.Range("N65").Value = -384864

$ws.Range("H66").Value = 96000
$ws.Range("I66").Value = 70000
$ws.Range("J66").Value = 122000
$ws.Range("K66").Value = 210000
$ws.Range("L66").Value = 366000
$ws.Range("M66").Value = -206568
$ws.Range("N66").Value = -372864

$ws.Range("H86").Value = 1797.2106
$ws.Range("I86").Value = 1988.5385
$ws.Range("K86").Value = 1988.5385
$ws.Range("M86").Value = -865.5385000000001

$ws.Range("H89").Value = 1797.2106
$ws.Range("I89").Value = 1988.5385
$ws.Range("K89").Value = 9942.692500000001
$ws.Range("M89").Value = -4326.692500000001

$ws.Range("H99").Value = 2370.7917
$ws.Range("I99").Value = 1778.8235
$ws.Range("K99").Value = 1778.8235
$ws.Range("M99").Value = -280.8235

$ws.Range("H105").Value = 1627.1111
$ws.Range("I105").Value = 1518
$ws.Range("K105").Value = 1518
$ws.Range("M105").Value = 229


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 150
$ws.Range("I22").Value = 150
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = -850

$ws.Range("H58").Value = 1581.2667
$ws.Range("I58").Value = 1530.0714
$ws.Range("K58").Value = 1530.0714
$ws.Range("M58").Value = -1327.0714

$ws.Range("H103").Value = 30978.445
$ws.Range("I103").Value = 23983
$ws.Range("J103").Value = 44969.332
$ws.Range("K103").Value = 23983
$ws.Range("L103").Value = 44969.332
$ws.Range("M103").Value = -22811
$ws.Range("N103").Value = -47313.332

$ws.Range("H111").Value = 59994.5
$ws.Range("J111").Value = 59994.5
$ws.Range("L111").Value = 59994.5
$ws.Range("N111").Value = -68174.5

$ws.Range("H132").Value = 2119.0667
$ws.Range("I132").Value = 1913.2858
$ws.Range("K132").Value = 5739.857400000001
$ws.Range("M132").Value = -3209.857400000001

$ws.Range("H136").Value = 1581.2667
$ws.Range("I136").Value = 1530.0714
$ws.Range("K136").Value = 4590.2142
$ws.Range("M136").Value = -2040.2142


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 31103.54
$ws.Range("I2").Value = 432
$ws.Range("J2").Value = 80178
$ws.Range("K2").Value = 2592
$ws.Range("L2").Value = 481068
$ws.Range("M2").Value = -2479
$ws.Range("N2").Value = -481294

$ws.Range("H38").Value = 58.68421
$ws.Range("J38").Value = 55.5
$ws.Range("L38").Value = 166.5
$ws.Range("N38").Value = -860.5


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 116.78571
$ws.Range("I2").Value = 68.5
$ws.Range("K2").Value = 68.5
$ws.Range("M2").Value = 44.5

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("M73").ClearContents()

$ws.Range("H97").Value = 1393.8636
$ws.Range("I97").Value = 1431.4375
$ws.Range("J97").Value = 1293.6666
$ws.Range("K97").Value = 1431.4375
$ws.Range("L97").Value = 1293.6666
$ws.Range("M97").Value = -935.4375
$ws.Range("N97").Value = -2285.6666


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 57884.633
$ws.Range("J7").Value = 96865.55
$ws.Range("L7").Value = 96865.55
$ws.Range("N7").Value = -97089.55

$ws.Range("H16").Value = 1353.4667
$ws.Range("I16").Value = 1538.4166
$ws.Range("K16").Value = 1538.4166
$ws.Range("M16").Value = -1368.4166

$ws.Range("H125").Value = 132000
$ws.Range("J125").Value = 132000
$ws.Range("L125").Value = 132000
$ws.Range("N125").Value = -141840

$ws.Range("H126").Value = 57884.633
$ws.Range("J126").Value = 96865.55
$ws.Range("L126").Value = 290596.65
$ws.Range("N126").Value = -295536.65

$ws.Range("H132").Value = 326244.03
$ws.Range("I132").Value = 590757.5600000001
$ws.Range("J132").Value = 61730.47
$ws.Range("K132").Value = 1772272.68
$ws.Range("L132").Value = 185191.41
$ws.Range("M132").Value = -1769742.68
$ws.Range("N132").Value = -190251.41

$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 89325
$ws.Range("J103").Value = 89325
$ws.Range("L103").Value = 89325
$ws.Range("N103").Value = -91669

$ws.Range("H125").Value = 73426.86
$ws.Range("J125").Value = 73426.86
$ws.Range("L125").Value = 73426.86
$ws.Range("N125").Value = -83266.86

$ws.Range("H132").Value = 3894.6287
$ws.Range("I132").Value = 3461.24
$ws.Range("K132").Value = 10383.72
$ws.Range("M132").Value = -7853.719999999999

